# Normalize the DATA_TYPE column values in the T_EMP and T_DEPT metadata
# sheets so they use generic/lower-case type names instead of the
# Oracle-specific ones (e.g. NUMBER(22) -> int, VARCHAR2(50) -> varchar(50),
# DATE -> date).

$wb = $excel.ActiveWorkbook

$wsEmp = $wb.Worksheets.Item("T_EMP")
$wsDept = $wb.Worksheets.Item("T_DEPT")

# Map of old Oracle-style data type text -> new generic data type text.
$typeMap = @{
    "NUMBER(22)"    = "int"
    "VARCHAR2(50)"  = "varchar(50)"
    "VARCHAR2(20)"  = "varchar(20)"
    "VARCHAR2(30)"  = "varchar(30)"
    "DATE"          = "date"
    "VARCHAR2(10)"  = "varchar(10)"
    "VARCHAR2(14)"  = "varchar(14)"
    "VARCHAR2(13)"  = "varchar(13)"
}

foreach ($ws in @($wsEmp, $wsDept)) {
    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        $cell = $ws.Cells.Item($r, 2)
        $current = $cell.Value2
        if ($typeMap.ContainsKey($current)) {
            $cell.Value = $typeMap[$current]
        }
    }
}
